$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 5: fill in D5:J5 with numeric values (TakeMinCount, PriceMin, PriceMax, ...)
$ws.Range("D5").Value = 0
$ws.Range("E5").Value = 100
$ws.Range("F5").Value = 10000
$ws.Range("G5").Value = 0
$ws.Range("H5").Value = 0
$ws.Range("I5").Value = 0
$ws.Range("J5").Value = 0

# Row 6: replace text value in D6 ("캐시재화") with numeric 0, and fill E6:J6
$ws.Range("D6").Value = 0
$ws.Range("E6").Value = 100
$ws.Range("F6").Value = 10000
$ws.Range("G6").Value = 0
$ws.Range("H6").Value = 0
$ws.Range("I6").Value = 0
$ws.Range("J6").Value = 0

# Row 7: replace text value in D7 ("게임 내에서 구할 수 있는 골드") with numeric 0, and fill E7:J7
$ws.Range("D7").Value = 0
$ws.Range("E7").Value = 100
$ws.Range("F7").Value = 10000
$ws.Range("G7").Value = 0
$ws.Range("H7").Value = 0
$ws.Range("I7").Value = 0
$ws.Range("J7").Value = 0

# Update sheet selection to A3:J19 with active cell A3
$ws.Activate()
$ws.Range("A3:J19").Select()

$wb.Save()
